$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the condition cell value: add a space after the colon
$ws.Range("B8").Value = "`$patient: Patient"

# Update the selection to match the new active cell / selected range
$ws.Range("B8").Select()
